# Update the PAY-GRADE-ELEMENTS sheet of the pay-grade-add CSV schema
# workbook: rename/expand the "salaryAmount" row into "joiningBonus", add
# several new optional fields, and keep "currency" plus a new
# "salaryAmount" row further down, ending with a new "countryCode" row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PAY-GRADE-ELEMENTS")

$numericDesc = "Column name is just a sample and can be changed to match your intake form variable. For numeric variables, provide comma-separated values in format: min,max,target (e.g., '100,200,150' where min=100, max=200, target=150). Min and max values are required for numeric fields; target is optional."
$genericDesc = "Column name is just a sample and can be changed. Refer to documentation."

# --- Row 3 currently holds "salaryAmount" (Numeric/Optional). Turn it
# into the new "joiningBonus" row, keeping the Numeric/Optional/blank/blank
# layout and just swapping the field name + description.
$ws.Range("A3").Value = "joiningBonus"
$ws.Range("B3").Value = $numericDesc

# --- Insert two new blank rows (formatted like the surrounding rows)
# right before the existing "currency" row (row 4), for "bonusPayFrequency"
# and "basePayFrequency". This pushes "currency" down to row 6.
$ws.Range("A4:F5").Insert()

$ws.Range("A4").Value = "bonusPayFrequency"
$ws.Range("B4").Value = $genericDesc
$ws.Range("C4").Value = "String"
$ws.Range("D4").Value = "Optional"

$ws.Range("A5").Value = "basePayFrequency"
$ws.Range("B5").Value = $genericDesc
$ws.Range("C5").Value = "String"
$ws.Range("D5").Value = "Optional"

# "currency" (now row 6) keeps its original values untouched.

# --- Insert six new blank rows after "currency" (row 6) for the new
# numeric compensation fields plus a trailing "countryCode" row.
$ws.Range("A7:F12").Insert()

$ws.Range("A7").Value = "salaryAmount"
$ws.Range("B7").Value = $numericDesc
$ws.Range("C7").Value = "Numeric"
$ws.Range("D7").Value = "Optional"

$ws.Range("A8").Value = "variableBonus"
$ws.Range("B8").Value = $numericDesc
$ws.Range("C8").Value = "Numeric"
$ws.Range("D8").Value = "Optional"

$ws.Range("A9").Value = "stockOptionsAmount"
$ws.Range("B9").Value = $numericDesc
$ws.Range("C9").Value = "Numeric"
$ws.Range("D9").Value = "Optional"

$ws.Range("A10").Value = "relocationBonus"
$ws.Range("B10").Value = $numericDesc
$ws.Range("C10").Value = "Numeric"
$ws.Range("D10").Value = "Optional"

$ws.Range("A11").Value = "totalCompensation"
$ws.Range("B11").Value = $numericDesc
$ws.Range("C11").Value = "Numeric"
$ws.Range("D11").Value = "Optional"

$ws.Range("A12").Value = "countryCode"
$ws.Range("B12").Value = $genericDesc
$ws.Range("C12").Value = "String"
$ws.Range("D12").Value = "Optional"
